$d = $word.ActiveDocument

# --- Locate the final paragraph ("我又吃了中饭", which carries a pPr/rFonts
# eastAsia hint) that sits right before the closing sectPr. Two brand new
# paragraphs need to land just above it:
#   1) "我又吃了中饭"              -- single run, rFonts hint=eastAsia
#   2) "2023年3月3日" + "星期五"   -- two runs, only the 2nd hinted eastAsia
# and then the original paragraph's text becomes "我再一次吃了中饭".

$lastIndex = $d.Paragraphs.Count
$targetPara = $d.Paragraphs($lastIndex)

# Insert an empty paragraph right before it, then fill it via InsertXML so we
# get full control of the run/rPr shape (avoids inheriting a stray pPr from
# the paragraph we split off of).
$null = $targetPara.Range.InsertParagraphBefore()

$lastIndex = $d.Paragraphs.Count
$newPara1 = $d.Paragraphs($lastIndex - 1)
$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>我又吃了中饭</w:t></w:r>' +
        '</w:p>'
$newPara1.Range.InsertXML($xml1)

# Insert the second new paragraph (date + weekday) right after the one we
# just filled in.
$null = $newPara1.Range.InsertParagraphAfter()

$lastIndex = $d.Paragraphs.Count
$newPara2 = $d.Paragraphs($lastIndex - 1)
$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:r><w:t>2023年3月3日</w:t></w:r>' +
        '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>星期五</w:t></w:r>' +
        '</w:p>'
$newPara2.Range.InsertXML($xml2)

# Finally, update the original closing paragraph's text.
$lastIndex = $d.Paragraphs.Count
$finalPara = $d.Paragraphs($lastIndex)
$finalPara.Range.Text = "我再一次吃了中饭"
